# "Most new technologies can start from 2019"
# Update the technology-start-year cells (previously 2020) to 2019 across
# the four worksheets of the workbook, and restore the view/selection state
# that Excel recorded for the RSD and RSD_OTH sheets while doing so.

$wb = $excel.ActiveWorkbook

# --- RSD --------------------------------------------------------------
$wsRSD = $wb.Worksheets.Item("RSD")
foreach ($addr in @("M7","M8","M11","M12","M13","M14")) {
    $wsRSD.Range($addr).Value = 2019
}

# --- RSD_buildings ------------------------------------------------------
$wsBuildings = $wb.Worksheets.Item("RSD_buildings")
foreach ($addr in @("K16","K21","K26")) {
    $wsBuildings.Range($addr).Value = 2019
}

# --- RSD_Boilers ----------------------------------------------------------
$wsBoilers = $wb.Worksheets.Item("RSD_Boilers")
foreach ($addr in @("I16","I17","I18","I19","I20","I21")) {
    $wsBoilers.Range($addr).Value = 2019
}

# --- RSD_OTH --------------------------------------------------------------
$wsOth = $wb.Worksheets.Item("RSD_OTH")
foreach ($addr in @("I14","I15","I16","I17","I18","I19")) {
    $wsOth.Range($addr).Value = 2019
}

# Update the recorded selection on RSD_OTH (was F25, now I14) before moving
# focus away from it.
$wsOth.Activate()
$wsOth.Range("I14").Select()

# Finally, RSD becomes the active sheet/tab, with M7 selected and the view
# scrolled so column C is at the left edge.
$wsRSD.Activate()
$wsRSD.Range("M7").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
